# Auto-generated Excel COM-interop script
# Applies updated cryptocurrency price (column D) and 1h volume change (column E)
# values to rows 2-51 of the active worksheet, matching the commit
# "Updated cryptos list on Thu Feb  8 17:32:49 UTC 2024 with GitHub Actions".
#
# Column D prices are stored as plain text in the original workbook (e.g.
# "44.958.50", "0.999", "1.00"). Excel's COM layer auto-converts numeric-
# looking strings to real numbers on assignment, which would both change the
# stored value (e.g. "1.00" -> 1) and lose the original text representation.
# To avoid that, each D cell is briefly switched to the Text number format
# ("@") before its new value is written, then restored to the workbook's
# default "Normal" style so no visible/serialized formatting change remains
# on the cell (matching the source diff, which only touches cell text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.080.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.429.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.49%  "
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.802.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.417.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.837"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.967.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "49.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.54%  "
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0763"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.24%  "
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.934.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.77%  "
$ws.Range("E50").Value = "  +5.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.80%  "
